# Update the Markov transition probability matrix on Sheet1 (rows 2-19,
# columns B-S) with refreshed values reflecting additional simulated games
# and updated optimization logic.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1939799331103679
$ws.Range("C2").Value = 0.5886287625418061
$ws.Range("J2").Value = 0.01672240802675585
$ws.Range("P2").Value = 0.1304347826086956
$ws.Range("S2").Value = 0.07023411371237458
$ws.Range("C3").Value = 0.01104972375690608
$ws.Range("J3").Value = 0.05524861878453038
$ws.Range("P3").Value = 0.7679558011049724
$ws.Range("S3").Value = 0.1657458563535912
$ws.Range("J4").Value = 0.04166666666666666
$ws.Range("O4").Value = 0.02083333333333333
$ws.Range("P4").Value = 0.7291666666666666
$ws.Range("S4").Value = 0.2083333333333333
$ws.Range("B6").Value = 0.07614213197969544
$ws.Range("D6").Value = 0.01015228426395939
$ws.Range("F6").Value = 0.06091370558375635
$ws.Range("J6").Value = 0.3248730964467005
$ws.Range("O6").Value = 0.01015228426395939
$ws.Range("Q6").Value = 0.1472081218274112
$ws.Range("R6").Value = 0.06598984771573604
$ws.Range("S6").Value = 0.3045685279187818
$ws.Range("B7").Value = 0.1101321585903084
$ws.Range("D7").Value = 0.03083700440528634
$ws.Range("E7").Value = 0.004405286343612335
$ws.Range("F7").Value = 0.03083700440528634
$ws.Range("J7").Value = 0.1718061674008811
$ws.Range("O7").Value = 0.01762114537444934
$ws.Range("Q7").Value = 0.1629955947136564
$ws.Range("R7").Value = 0.06607929515418502
$ws.Range("S7").Value = 0.4052863436123348
$ws.Range("B8").Value = 0.08029197080291971
$ws.Range("D8").Value = 0.0218978102189781
$ws.Range("F8").Value = 0.08029197080291971
$ws.Range("J8").Value = 0.1386861313868613
$ws.Range("O8").Value = 0.0170316301703163
$ws.Range("Q8").Value = 0.2214111922141119
$ws.Range("R8").Value = 0.07785888077858881
$ws.Range("S8").Value = 0.3625304136253041
$ws.Range("B9").Value = 0.08333333333333333
$ws.Range("D9").Value = 0.01666666666666667
$ws.Range("F9").Value = 0.05833333333333333
$ws.Range("J9").Value = 0.1166666666666667
$ws.Range("O9").Value = 0.03333333333333333
$ws.Range("Q9").Value = 0.1833333333333333
$ws.Range("R9").Value = 0.1
$ws.Range("S9").Value = 0.4083333333333333
$ws.Range("B10").Value = 0.1183800623052959
$ws.Range("D10").Value = 0.02336448598130841
$ws.Range("F10").Value = 0.05218068535825545
$ws.Range("J10").Value = 0.1355140186915888
$ws.Range("O10").Value = 0.01713395638629283
$ws.Range("Q10").Value = 0.2367601246105919
$ws.Range("R10").Value = 0.05841121495327103
$ws.Range("S10").Value = 0.3582554517133956
$ws.Range("G11").Value = 0.1581920903954802
$ws.Range("J11").Value = 0.09887005649717515
$ws.Range("K11").Value = 0.2062146892655367
$ws.Range("L11").Value = 0.5282485875706214
$ws.Range("S11").Value = 0.008474576271186441
$ws.Range("G12").Value = 0.7225130890052356
$ws.Range("J12").Value = 0.225130890052356
$ws.Range("K12").Value = 0.01047120418848168
$ws.Range("L12").Value = 0.01047120418848168
$ws.Range("S12").Value = 0.03141361256544502
$ws.Range("G13").Value = 0.6530612244897959
$ws.Range("J13").Value = 0.3265306122448979
$ws.Range("S13").Value = 0.02040816326530612
$ws.Range("F15").Value = 0.02873563218390805
$ws.Range("H15").Value = 0.1264367816091954
$ws.Range("I15").Value = 0.07471264367816093
$ws.Range("J15").Value = 0.4022988505747127
$ws.Range("K15").Value = 0.07471264367816093
$ws.Range("M15").Value = 0.02873563218390805
$ws.Range("O15").Value = 0.04022988505747126
$ws.Range("S15").Value = 0.2241379310344828
$ws.Range("F16").Value = 0.01442307692307692
$ws.Range("H16").Value = 0.1538461538461539
$ws.Range("I16").Value = 0.05288461538461538
$ws.Range("J16").Value = 0.4519230769230769
$ws.Range("K16").Value = 0.1586538461538461
$ws.Range("M16").Value = 0.01442307692307692
$ws.Range("O16").Value = 0.02884615384615385
$ws.Range("S16").Value = 0.125
$ws.Range("F17").Value = 0.02049180327868852
$ws.Range("H17").Value = 0.1864754098360656
$ws.Range("I17").Value = 0.05737704918032787
$ws.Range("J17").Value = 0.4118852459016393
$ws.Range("K17").Value = 0.1413934426229508
$ws.Range("M17").Value = 0.01229508196721311
$ws.Range("N17").Value = 0.002049180327868853
$ws.Range("O17").Value = 0.06352459016393443
$ws.Range("S17").Value = 0.1045081967213115
$ws.Range("F18").Value = 0.02083333333333333
$ws.Range("H18").Value = 0.1736111111111111
$ws.Range("I18").Value = 0.04166666666666666
$ws.Range("J18").Value = 0.5138888888888888
$ws.Range("K18").Value = 0.09027777777777778
$ws.Range("M18").Value = 0.006944444444444444
$ws.Range("O18").Value = 0.03472222222222222
$ws.Range("S18").Value = 0.1180555555555556
$ws.Range("F19").Value = 0.02327663384064458
$ws.Range("H19").Value = 0.2175470008952551
$ws.Range("I19").Value = 0.05640107430617726
$ws.Range("J19").Value = 0.3724261414503133
$ws.Range("K19").Value = 0.1316025067144136
$ws.Range("M19").Value = 0.03043867502238138
$ws.Range("O19").Value = 0.05729632945389436
$ws.Range("S19").Value = 0.1110116383169203

Write-Output "Updated 108 cells"
